$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.056.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.83%  '
$ws.Range("D3").Value = "'1.693.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'220.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.44%  '
$ws.Range("D6").Value = "'0.534"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.45%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'29.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.95%  '
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("E10").Value = '  +4.87%  '
$ws.Range("D11").Value = "'0.0910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("D12").Value = "'1.933.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.45%  '
$ws.Range("D13").Value = "'1.686.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").Value = "'10.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.63%  '
$ws.Range("D15").Value = "'0.609"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("E16").Value = '  +6.83%  '
$ws.Range("D17").Value = "'31.063.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'67.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("D19").Value = "'248.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.07%  '
$ws.Range("D20").Value = "'0.0₃0720"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = "'4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.86%  '
$ws.Range("D23").Value = "'10.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").Value = "'158.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = "'15.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.59%  '
$ws.Range("D27").Value = "'0.112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("D28").Value = "'6.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = "'0.0501"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("D31").Value = "'3.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.81%  '
$ws.Range("E32").Value = '  +3.86%  '
$ws.Range("D33").Value = "'3.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.65%  '
$ws.Range("D34").Value = "'1.514.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = '  +2.43%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = "'83.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.12%  '
$ws.Range("D38").Value = "'0.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.16%  '
$ws.Range("E39").Value = '  +4.53%  '
$ws.Range("D40").Value = "'2.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.09%  '
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").Value = "'2.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.35%  '
$ws.Range("D43").Value = "'0.847"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("D44").Value = "'0.0505"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").Value = "'1.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.94%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = "'52.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.83%  '
$ws.Range("D48").Value = "'5.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.09%  '
$ws.Range("D49").Value = "'1.822.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("E50").Value = '  +9.65%  '
$ws.Range("D51").Value = "'93.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.25%  '
